# Update cryptocurrency price/volume figures per the Oct 4 2024 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.342.56"
$ws.Range("E2").Value = "'  +2.20%  "
$ws.Range("D3").Value = "'2.425.08"
$ws.Range("E3").Value = "'  +3.09%  "
$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("D5").Value = "'556.70"
$ws.Range("E5").Value = "'  +2.02%  "
$ws.Range("D6").Value = "'143.83"
$ws.Range("E6").Value = "'  +4.82%  "
$ws.Range("E7").Value = "'  -0.04%  "
$ws.Range("E8").Value = "'  +1.71%  "
$ws.Range("D9").Value = "'2.424.96"
$ws.Range("E9").Value = "'  +3.05%  "
$ws.Range("E10").Value = "'  +5.14%  "
$ws.Range("E11").Value = "'  +0.15%  "
$ws.Range("E12").Value = "'  +1.72%  "
$ws.Range("D13").Value = "'0.353"
$ws.Range("E13").Value = "'  +3.13%  "
$ws.Range("D14").Value = "'26.29"
$ws.Range("E14").Value = "'  +6.32%  "
$ws.Range("E15").Value = "'  +9.40%  "
$ws.Range("D16").Value = "'2.864.21"
$ws.Range("E16").Value = "'  +3.14%  "
$ws.Range("D17").Value = "'62.275.69"
$ws.Range("E17").Value = "'  +2.32%  "
$ws.Range("D18").Value = "'2.430.16"
$ws.Range("E18").Value = "'  +3.42%  "
$ws.Range("D19").Value = "'11.12"
$ws.Range("E19").Value = "'  +4.72%  "
$ws.Range("D20").Value = "'325.40"
$ws.Range("E20").Value = "'  +1.65%  "
$ws.Range("D21").Value = "'4.18"
$ws.Range("E21").Value = "'  +1.52%  "
$ws.Range("E22").Value = "'  +3.32%  "
$ws.Range("E23").Value = "'  +0.20%  "
$ws.Range("D24").Value = "'1.79"
$ws.Range("E24").Value = "'  +2.83%  "
$ws.Range("D25").Value = "'65.06"
$ws.Range("E25").Value = "'  +2.80%  "
$ws.Range("D26").Value = "'9.18"
$ws.Range("E26").Value = "'  +9.92%  "
$ws.Range("D27").Value = "'574.67"
$ws.Range("E27").Value = "'  +14.55%  "
$ws.Range("D28").Value = "'0.0₃0955"
$ws.Range("E28").Value = "'  +9.65%  "
$ws.Range("D29").Value = "'2.545.14"
$ws.Range("E29").Value = "'  +3.17%  "
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "'  -0.02%  "
$ws.Range("D31").Value = "'8.42"
$ws.Range("E32").Value = "'  +5.60%  "
$ws.Range("E33").Value = "'  +1.61%  "
$ws.Range("E34").Value = "'  +4.20%  "
$ws.Range("E35").Value = "'  +4.33%  "
$ws.Range("E36").Value = "'  +8.58%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "'  +0.04%  "
$ws.Range("D38").Value = "'4.83"
$ws.Range("E38").Value = "'  +4.34%  "
$ws.Range("E39").Value = "'  +2.19%  "
$ws.Range("E40").Value = "'  +3.83%  "
$ws.Range("E41").Value = "'  +1.66%  "
$ws.Range("D42").Value = "'148.39"
$ws.Range("E42").Value = "'  +4.17%  "
$ws.Range("E43").Value = "'  +0.02%  "
$ws.Range("D44").Value = "'41.70"
$ws.Range("E44").Value = "'  +2.84%  "
$ws.Range("D45").Value = "'2.32"
$ws.Range("E45").Value = "'  +12.04%  "
$ws.Range("D46").Value = "'152.35"
$ws.Range("E46").Value = "'  +6.73%  "
$ws.Range("E47").Value = "'  +2.05%  "
$ws.Range("D48").Value = "'0.0546"
$ws.Range("E48").Value = "'  +5.36%  "
$ws.Range("D49").Value = "'20.48"
$ws.Range("E49").Value = "'  +6.88%  "
$ws.Range("E50").Value = "'  +3.73%  "
$ws.Range("D51").Value = "'0.0229"
